# [HOTFIX] 214893 Entitlement export add national id
# Adds a new "national_id" column (M) to the payment list export sheet,
# mirroring the formatting of the existing "status" column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column L's formatting (fill/border/number format) onto the new
# column M first (rows 1-10 only, matching the sheet's used range), so the
# subsequent value writes keep that style.
$ws.Range("L1:L10").Copy() | Out-Null
$ws.Range("M1:M10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New header cell
$ws.Range("M1").Value = "national_id"

# Row 2 gets a sample national id ("123") - keep it as text like the rest
# of the id-like columns in this sheet (e.g. payment_id "123123").
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "123"

# Row 3 has no national_id value, cell stays empty (style only).
$ws.Range("M3").ClearContents() | Out-Null

# Column widths: D:E grow slightly (23.4531 -> 23.5) and the new column M
# joins K:L's wide "36.6719" sizing.
$ws.Columns.Item(4).ColumnWidth = 22.666666666666668
$ws.Columns.Item(5).ColumnWidth = 22.666666666666668
$ws.Columns.Item(13).ColumnWidth = 35.838566666666665
